# This script applies a weekly update to the "Chirimoya" price sheet:
# it inserts 4 new rows (a new date block, serial 45120 / 2023-07-13) right
# before the existing 44518 block (old row 222), shifting everything that
# used to occupy rows 222:277 down to rows 226:281.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows at 222, pushing the old 222:277 block down to 226:281.
$ws.Rows("222:225").Insert()

# Common (fixed) values shared by every data row in this table.
$marketId   = 9
$market     = "Vega Central Mapocho de Santiago"
$region     = "Metropolitana"
$codreg     = 13
$tipo       = "Fruta"
$productoId = 100107
$producto   = "Otros"
$categoriaId = 100107002
$categoria  = "Chirimoya"
$variedad   = "Cultivar IV Región"

# New date-block data (fecha serial 45120), quality rows: Especial, Extra
# (doble especial), Primera, Segunda - origin "Provincia del Elquí".
$newRows = @(
    @{ Row=222; Calidad="Especial";                Volumen=200; Min=21600; Max=21600; Prom=21600; Unidad="`$/bandeja 8 kilos"; Origen="Provincia del Elquí"; PrecioKg=2700; KgUnidad=8 },
    @{ Row=223; Calidad="Extra (doble especial)";   Volumen=120; Min=24000; Max=24000; Prom=24000; Unidad="`$/bandeja 8 kilos"; Origen="Provincia del Elquí"; PrecioKg=3000; KgUnidad=8 },
    @{ Row=224; Calidad="Primera";                  Volumen=100; Min=17600; Max=17600; Prom=17600; Unidad="`$/bandeja 8 kilos"; Origen="Provincia del Elquí"; PrecioKg=2200; KgUnidad=8 },
    @{ Row=225; Calidad="Segunda";                  Volumen=200; Min=14400; Max=14400; Prom=14400; Unidad="`$/bandeja 8 kilos"; Origen="Provincia del Elquí"; PrecioKg=1800; KgUnidad=8 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $marketId
    $ws.Cells.Item($row, 2).Value  = $market
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = 45120
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
